# Applies cell-value updates (updated odds/stats) to Sheet1 of the FlashScore workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19
$ws.Cells.Item(19, 7).Value = 3.7  # G19: 3.6 -> 3.7
$ws.Cells.Item(19, 8).Value = 2.6  # H19: 2.63 -> 2.6
$ws.Cells.Item(19, 10).Value = 4.75  # J19: 4.5 -> 4.75
$ws.Cells.Item(19, 13).Value = 1.2  # M19: 1.18 -> 1.2
$ws.Cells.Item(19, 14).Value = 4.33  # N19: 4.5 -> 4.33
$ws.Cells.Item(19, 15).Value = 1.83  # O19: 1.73 -> 1.83
$ws.Cells.Item(19, 16).Value = 1.83  # P19: 2 -> 1.83
$ws.Cells.Item(19, 17).Value = 4  # Q19: 3.6 -> 4
$ws.Cells.Item(19, 18).Value = 1.25  # R19: 1.29 -> 1.25
$ws.Cells.Item(19, 19).Value = 1.83  # S19: 1.8 -> 1.83
$ws.Cells.Item(19, 20).Value = 1.98  # T19: 2 -> 1.98
$ws.Cells.Item(19, 23).Value = 6.5  # W19: 7 -> 6.5
$ws.Cells.Item(19, 25).Value = 17  # Y19: 15 -> 17
$ws.Cells.Item(19, 28).Value = 67  # AB19: 51 -> 67
$ws.Cells.Item(19, 29).Value = 4.33  # AC19: 4.5 -> 4.33
$ws.Cells.Item(19, 30).Value = 6  # AD19: 5.5 -> 6
$ws.Cells.Item(19, 31).Value = 26  # AE19: 23 -> 26
$ws.Cells.Item(19, 32).Value = 126  # AF19: 101 -> 126
$ws.Cells.Item(19, 38).Value = 34  # AL19: 29 -> 34
$ws.Cells.Item(19, 41).Value = 26  # AO19: 23 -> 26
$ws.Cells.Item(19, 44).Value = 201  # AR19: 151 -> 201
$ws.Cells.Item(19, 47).Value = 12  # AU19: 11 -> 12
$ws.Cells.Item(19, 48).Value = 126  # AV19: 101 -> 126
$ws.Cells.Item(19, 50).Value = 19  # AX19: 17 -> 19
$ws.Cells.Item(19, 52).Value = 67  # AZ19: 51 -> 67

# Row 20
$ws.Cells.Item(20, 7).Value = 3.7  # G20: 3.6 -> 3.7
$ws.Cells.Item(20, 8).Value = 2.9  # H20: 2.88 -> 2.9
$ws.Cells.Item(20, 9).Value = 2.3  # I20: 2.35 -> 2.3
$ws.Cells.Item(20, 12).Value = 3.25  # L20: 3.4 -> 3.25
$ws.Cells.Item(20, 14).Value = 5  # N20: 4.75 -> 5
$ws.Cells.Item(20, 24).Value = 17  # X20: 15 -> 17
$ws.Cells.Item(20, 28).Value = 67  # AB20: 51 -> 67
$ws.Cells.Item(20, 34).Value = 4.75  # AH20: 5 -> 4.75
$ws.Cells.Item(20, 35).Value = 8.5  # AI20: 9 -> 8.5
$ws.Cells.Item(20, 37).Value = 21  # AK20: 23 -> 21
$ws.Cells.Item(20, 40).Value = 5.5  # AN20: 5 -> 5.5
$ws.Cells.Item(20, 50).Value = 15  # AX20: 17 -> 15

# Row 21
$ws.Cells.Item(21, 13).Value = 1.14  # M21: 1.13 -> 1.14
$ws.Cells.Item(21, 14).Value = 5.5  # N21: 6 -> 5.5
$ws.Cells.Item(21, 15).Value = 1.67  # O21: 1.62 -> 1.67
$ws.Cells.Item(21, 16).Value = 2.1  # P21: 2.2 -> 2.1

# Row 34
$ws.Cells.Item(34, 17).Value = 1.73  # Q34: 1.75 -> 1.73
$ws.Cells.Item(34, 18).Value = 2.08  # R34: 2.05 -> 2.08

# Row 35
$ws.Cells.Item(35, 13).Value = 1.03  # M35: 1.04 -> 1.03
$ws.Cells.Item(35, 14).Value = 15  # N35: 13 -> 15

# Row 36
$ws.Cells.Item(36, 10).Value = 2.3  # J36: 2.25 -> 2.3
$ws.Cells.Item(36, 33).Value = 501  # AG36: 451 -> 501

# Row 37
$ws.Cells.Item(37, 7).Value = 3.2  # G37: 3.1 -> 3.2
$ws.Cells.Item(37, 9).Value = 2.38  # I37: 2.4 -> 2.38

# Row 38
$ws.Cells.Item(38, 13).Value = 1.11  # M38: 1.1 -> 1.11
$ws.Cells.Item(38, 14).Value = 6.5  # N38: 7 -> 6.5
$ws.Cells.Item(38, 15).Value = 1.5  # O38: 1.44 -> 1.5
$ws.Cells.Item(38, 16).Value = 2.63  # P38: 2.75 -> 2.63
$ws.Cells.Item(38, 17).Value = 2.5  # Q38: 2.4 -> 2.5
$ws.Cells.Item(38, 18).Value = 1.5  # R38: 1.53 -> 1.5

# Row 39
$ws.Cells.Item(39, 13).Value = 1.08  # M39: 1.07 -> 1.08
$ws.Cells.Item(39, 14).Value = 8  # N39: 9 -> 8
$ws.Cells.Item(39, 17).Value = 2.15  # Q39: 2.1 -> 2.15
$ws.Cells.Item(39, 18).Value = 1.67  # R39: 1.7 -> 1.67

# Row 41
$ws.Cells.Item(41, 13).Value = 1.11  # M41: 1.1 -> 1.11
$ws.Cells.Item(41, 14).Value = 6.5  # N41: 7 -> 6.5
$ws.Cells.Item(41, 15).Value = 1.53  # O41: 1.5 -> 1.53
$ws.Cells.Item(41, 16).Value = 2.38  # P41: 2.5 -> 2.38
$ws.Cells.Item(41, 17).Value = 2.7  # Q41: 2.6 -> 2.7
$ws.Cells.Item(41, 18).Value = 1.44  # R41: 1.48 -> 1.44
$ws.Cells.Item(41, 19).Value = 1.62  # S41: 1.57 -> 1.62
$ws.Cells.Item(41, 20).Value = 2.2  # T41: 2.25 -> 2.2
$ws.Cells.Item(41, 32).Value = 81  # AF41: 67 -> 81
$ws.Cells.Item(41, 34).Value = 6  # AH41: 6.5 -> 6
$ws.Cells.Item(41, 35).Value = 10  # AI41: 11 -> 10
$ws.Cells.Item(41, 38).Value = 26  # AL41: 23 -> 26
$ws.Cells.Item(41, 45).Value = 351  # AS41: 301 -> 351
$ws.Cells.Item(41, 46).Value = 2.2  # AT41: 2.25 -> 2.2
$ws.Cells.Item(41, 47).Value = 9.5  # AU41: 9 -> 9.5
$ws.Cells.Item(41, 51).Value = 34  # AY41: 29 -> 34
$ws.Cells.Item(41, 53).Value = 101  # BA41: 81 -> 101

# Row 86
$ws.Cells.Item(86, 16).Value = 3.72  # P86: 3.7 -> 3.72

# Row 87
$ws.Cells.Item(87, 7).Value = 2.32  # G87: 2.65 -> 2.32
$ws.Cells.Item(87, 8).Value = 3.15  # H87: 2.9 -> 3.15
$ws.Cells.Item(87, 9).Value = 2.77  # I87: 2.57 -> 2.77
$ws.Cells.Item(87, 10).Value = 2.8  # J87: 3.15 -> 2.8
$ws.Cells.Item(87, 11).Value = 2.12  # K87: 2.05 -> 2.12
$ws.Cells.Item(87, 12).Value = 3.3  # L87: 3.1 -> 3.3
$ws.Cells.Item(87, 13).Value = 1.03  # M87: 1.01 -> 1.03
$ws.Cells.Item(87, 14).Value = 11.9  # N87: 11 -> 11.9
$ws.Cells.Item(87, 19).Value = 1.28  # S87: 1.29 -> 1.28
$ws.Cells.Item(87, 21).Value = 1.47  # U87: 1.46 -> 1.47
$ws.Cells.Item(87, 22).Value = 2.6  # V87: 2.63 -> 2.6
$ws.Cells.Item(87, 23).Value = 10  # W87: 10.5 -> 10
$ws.Cells.Item(87, 24).Value = 13  # X87: 15 -> 13
$ws.Cells.Item(87, 25).Value = 7.7  # Y87: 8 -> 7.7
$ws.Cells.Item(87, 26).Value = 22  # Z87: 29 -> 22
$ws.Cells.Item(87, 27).Value = 13.5  # AA87: 16 -> 13.5
$ws.Cells.Item(87, 28).Value = 15  # AB87: 16 -> 15
$ws.Cells.Item(87, 29).Value = 13  # AC87: 12.5 -> 13
$ws.Cells.Item(87, 30).Value = 5.9  # AD87: 5.5 -> 5.9
$ws.Cells.Item(87, 31).Value = 8.75  # AE87: 7.9 -> 8.75
$ws.Cells.Item(87, 32).Value = 25  # AF87: 22 -> 25
$ws.Cells.Item(87, 33).Value = 120  # AG87: 100 -> 120
$ws.Cells.Item(87, 36).Value = 8.5  # AJ87: 8 -> 8.5
$ws.Cells.Item(87, 37).Value = 28  # AK87: 27 -> 28
$ws.Cells.Item(87, 38).Value = 17  # AL87: 15.5 -> 17
$ws.Cells.Item(87, 39).Value = 18  # AM87: 16 -> 18
$ws.Cells.Item(87, 40).Value = 4.6  # AN87: 5 -> 4.6
$ws.Cells.Item(87, 41).Value = 12  # AO87: 14.5 -> 12
$ws.Cells.Item(87, 42).Value = 15.5  # AP87: 16.5 -> 15.5
$ws.Cells.Item(87, 43).Value = 45  # AQ87: 60 -> 45
$ws.Cells.Item(87, 44).Value = 55  # AR87: 70 -> 55
$ws.Cells.Item(87, 46).Value = 3.05  # AT87: 2.95 -> 3.05
$ws.Cells.Item(87, 47).Value = 6  # AU87: 5.6 -> 6
$ws.Cells.Item(87, 48).Value = 40  # AV87: 35 -> 40
$ws.Cells.Item(87, 49).Value = 5.1  # AW87: 4.9 -> 5.1
$ws.Cells.Item(87, 50).Value = 15.5  # AX87: 14 -> 15.5
$ws.Cells.Item(87, 51).Value = 18.5  # AY87: 16.5 -> 18.5
$ws.Cells.Item(87, 52).Value = 65  # AZ87: 60 -> 65
$ws.Cells.Item(87, 53).Value = 80  # BA87: 70 -> 80
$ws.Cells.Item(87, 54).Value = 175  # BB87: 150 -> 175
$ws.Cells.Item(87, 55).Value = 500  # BC87: 450 -> 500

# Row 106
$ws.Cells.Item(106, 13).Value = 1.06  # M106: 1.08 -> 1.06
$ws.Cells.Item(106, 14).Value = 10  # N106: 8 -> 10
